$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Status value for row 11 from "Y4" to "Y7"
$ws.Range("N11").Value = "Y7"

# Update numeric values in row 11 to reflect the new SQ1 metrics
$ws.Range("A11").Value = 687
$ws.Range("B11").Value = 26
$ws.Range("C11").Value = 661
$ws.Range("D11").Value = 3.8
$ws.Range("E11").Value = 26
$ws.Range("F11").Value = 26
$ws.Range("I11").Value = 33.9
$ws.Range("J11").Value = 233
$ws.Range("M11").Value = 31.18
$ws.Range("Q11").Value = 26
$ws.Range("R11").Value = 3.8
$ws.Range("T11").Value = 33.9
$ws.Range("U11").Value = 31.18
